$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 155.8125
$ws.Range("I11").Value = 155.8125
$ws.Range("K11").Value = 155.8125
$ws.Range("M11").Value = -15.8125
$ws.Range("H12").Value = 380.73334
$ws.Range("I12").Value = 382.92856
$ws.Range("K12").Value = 382.92856
$ws.Range("M12").Value = -212.92856
$ws.Range("H28").Value = 707.4828
$ws.Range("I28").Value = 304.0909
$ws.Range("J28").Value = 1975.2858
$ws.Range("K28").Value = 304.0909
$ws.Range("L28").Value = 1975.2858
$ws.Range("M28").Value = 180.9091
$ws.Range("N28").Value = -2945.2858
$ws.Range("H38").Value = 1115.7142
$ws.Range("J38").Value = 7000
$ws.Range("L38").Value = 21000
$ws.Range("N38").Value = -21744
$ws.Range("H51").Value = 500000000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H53").Value = 62501050
$ws.Range("I53").Value = 71429690
$ws.Range("K53").Value = 71429690
$ws.Range("M53").Value = -71429053
$ws.Range("H58").Value = 451.25
$ws.Range("I58").Value = 451.25
$ws.Range("K58").Value = 1353.75
$ws.Range("M58").Value = -1203.75
$ws.Range("H92").Value = 2351.75
$ws.Range("I92").Value = 2534.6667
$ws.Range("K92").Value = 2534.6667
$ws.Range("M92").Value = -1286.6667
$ws.Range("H99").Value = 532.5
$ws.Range("J99").Value = 999.5
$ws.Range("L99").Value = 2998.5
$ws.Range("N99").Value = -5994.5
$ws.Range("H106").Value = 2512.5386
$ws.Range("I106").Value = 2648.2
$ws.Range("K106").Value = 2648.2
$ws.Range("M106").Value = -2017.2
$ws.Range("H137").Value = 5014837
$ws.Range("I137").Value = 25001450
$ws.Range("J137").Value = 18183.625
$ws.Range("K137").Value = 75004350
$ws.Range("L137").Value = 54550.875
$ws.Range("M137").Value = -75001800
$ws.Range("N137").Value = -59650.875
$ws.Range("H138").Value = 3444.0667
$ws.Range("J138").Value = 3032.0364
$ws.Range("L138").Value = 9096.109199999999
$ws.Range("N138").Value = -19376.1092

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1077.8269
$ws.Range("I32").Value = 1077.8269
$ws.Range("K32").Value = 1077.8269
$ws.Range("M32").Value = -790.8269
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 134314.08
$ws.Range("I74").Value = 151969.66
$ws.Range("K74").Value = 151969.66
$ws.Range("M74").Value = -151095.66
$ws.Range("H77").Value = 134314.08
$ws.Range("I77").Value = 151969.66
$ws.Range("K77").Value = 759848.3
$ws.Range("M77").Value = -755480.3
$ws.Range("H132").Value = 1392.8292
$ws.Range("I132").Value = 792.1539
$ws.Range("K132").Value = 2376.4617
$ws.Range("M132").Value = 153.5383000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2786.4583
$ws.Range("I134").Value = 1949.8182
$ws.Range("J134").Value = 3494.3845
$ws.Range("K134").Value = 5849.4546
$ws.Range("L134").Value = 10483.1535
$ws.Range("M134").Value = -3314.4546
$ws.Range("N134").Value = -15553.1535
$ws.Range("H138").Value = 64740.668
$ws.Range("J138").Value = 64740.668
$ws.Range("L138").Value = 64740.668
$ws.Range("N138").Value = -75020.66800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5525.8335
$ws.Range("I31").Value = 3637.5715
$ws.Range("J31").Value = 8169.4
$ws.Range("K31").Value = 3637.5715
$ws.Range("L31").Value = 8169.4
$ws.Range("M31").Value = -3342.5715
$ws.Range("N31").Value = -8759.4
$ws.Range("H34").Value = 5525.8335
$ws.Range("I34").Value = 3637.5715
$ws.Range("J34").Value = 8169.4
$ws.Range("K34").Value = 3637.5715
$ws.Range("L34").Value = 8169.4
$ws.Range("M34").Value = -3435.5715
$ws.Range("N34").Value = -8573.4
$ws.Range("H86").Value = 7799.5
$ws.Range("I86").Value = 7963.3335
$ws.Range("J86").Value = 7308
$ws.Range("K86").Value = 7963.3335
$ws.Range("L86").Value = 7308
$ws.Range("M86").Value = -6840.3335
$ws.Range("N86").Value = -9554
$ws.Range("H89").Value = 7799.5
$ws.Range("I89").Value = 7963.3335
$ws.Range("J89").Value = 7308
$ws.Range("K89").Value = 41000
$ws.Range("L89").Value = 36540
$ws.Range("M89").Value = -34200.6675
$ws.Range("N89").Value = -47772

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11906969
$ws.Range("I131").Value = 27779860
$ws.Range("J131").Value = 2301.25
$ws.Range("K131").Value = 83339580
$ws.Range("L131").Value = 6903.75
$ws.Range("M131").Value = -83334540
$ws.Range("N131").Value = -16983.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 17205850
$ws.Range("I11").Value = 40000780
$ws.Range("J11").Value = 4542000
$ws.Range("K11").Value = 40000780
$ws.Range("L11").Value = 4542000
$ws.Range("M11").Value = -40000641
$ws.Range("N11").Value = -4542278
$ws.Range("H70").Value = 31255270
$ws.Range("J70").Value = 6376.6665
$ws.Range("L70").Value = 6376.6665
$ws.Range("N70").Value = -6916.6665
$ws.Range("H73").Value = 31255270
$ws.Range("J73").Value = 6376.6665
$ws.Range("L73").Value = 6376.6665
$ws.Range("N73").Value = -8248.666499999999
$ws.Range("H132").Value = 2745.4
$ws.Range("I132").Value = 2681.75
$ws.Range("K132").Value = 8045.25
$ws.Range("M132").Value = -5515.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2700
$ws.Range("I16").Value = 833.3333
$ws.Range("J16").Value = 5500
$ws.Range("K16").Value = 833.3333
$ws.Range("L16").Value = 5500
$ws.Range("M16").Value = -663.3333
$ws.Range("N16").Value = -5840
$ws.Range("H136").Value = 5911.5713
$ws.Range("J136").Value = 7547.9
$ws.Range("L136").Value = 22643.7
$ws.Range("N136").Value = -27743.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 458520.75
$ws.Range("I136").Value = 3091.5
$ws.Range("J136").Value = 1597093.9
$ws.Range("K136").Value = 9274.5
$ws.Range("L136").Value = 4791281.699999999
$ws.Range("M136").Value = -6724.5
$ws.Range("N136").Value = -4796381.699999999
